# #327 Ajout des profils d'acces a58d18c1e8091c98efec92c8c093b361a253eee5
#
# 1. Update the "Date" metadata value.
# 2. Swap the two "Mapping" columns (AK <-> AL) on the "Elements" sheet,
#    including the header text, the data cells for rows 2-6, and the
#    (bestFit) column widths.

$wb = $excel.ActiveWorkbook

# --- 1. Metadata!B8 : Date ---------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-19T13:17:15+00:00"

# --- 2. Elements : swap column AK and column AL -------------------------
$elements = $wb.Worksheets.Item("Elements")

# Swap the header + data values for rows 1 (header) through 6.
# NB: reading back a freshly-written `.Value` is unreliable on this host,
# so use `.Value2` for reads (writes still go through `.Value`).
for ($row = 1; $row -le 6; $row++) {
    $akCell = $elements.Cells.Item($row, 37)
    $alCell = $elements.Cells.Item($row, 38)

    $akValue = $akCell.Value2
    $alValue = $alCell.Value2

    $akCell.Value = $alValue
    $alCell.Value = $akValue
}

# Swap the (bestFit) column widths to match the new content: AK becomes the
# wide "Spécification métier..." column, AL becomes the narrower "RIM
# Mapping" column.
$elements.Columns.Item(37).ColumnWidth = 73.333333333
$elements.Columns.Item(38).ColumnWidth = 24.166666667
